$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": insert a new date column at DS (14-nov) and shift the
# existing 01-oct..31-oct columns one slot to the right (DT..EX).
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Columns("DS:DS").Insert()

$wsPrix.Range("DS1").Value = "14-nov"
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 123).Value = "-"
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the new daily quote row (2025-11-12 / 29).
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDate = $wsGaz.Cells.Item(151, 1)
$gazDate.NumberFormat = "@"
$gazDate.Value = "2025-11-12"
$gazDate.Style = "Normal"
$wsGaz.Cells.Item(151, 2).Value = 29

# ---------------------------------------------------------------------------
# Sheet "CO2": append the new daily quote row (2025-11-12 / 81.75).
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$co2Date = $wsCO2.Cells.Item(151, 1)
$co2Date.NumberFormat = "@"
$co2Date.Value = "2025-11-12"
$co2Date.Style = "Normal"
$wsCO2.Cells.Item(151, 2).Value = 81.75
